$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 1080 (shifts all existing rows 1080..1161 down by one,
# matching Excel's default "insert copies formatting from the row above" so the
# date cell in column D keeps its date style automatically).
$ws.Rows.Item(1080).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A1080").Value = 9
$ws.Range("B1080").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1080").Value = "Metropolitana"
$ws.Range("D1080").Value = 45021
$ws.Range("E1080").Value = 13
$ws.Range("F1080").Value = "Fruta"
$ws.Range("G1080").Value = 100102
$ws.Range("H1080").Value = "Cítricos"
$ws.Range("I1080").Value = 100102005
$ws.Range("J1080").Value = "Naranja"
$ws.Range("K1080").Value = "Valencia"
$ws.Range("L1080").Value = "Primera"
$ws.Range("M1080").Value = 380
$ws.Range("N1080").Value = 12500
$ws.Range("O1080").Value = 12500
$ws.Range("P1080").Value = 12500
$ws.Range("Q1080").Value = "`$/caja 18 kilos granel"
$ws.Range("R1080").Value = "Provincia de Quillota"
$ws.Range("S1080").Value = 694
$ws.Range("T1080").Value = 18
